$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values: force text format so numeric-looking strings
# (e.g. "582.96") are not auto-converted to numbers by Excel, matching the
# original inline-string cell type.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D11", "D13", "D15", "D16", "D18", "D22", "D26", "D27", "D29", "D31", "D33", "D35", "D37", "D39", "D40", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.538.21"
$ws.Range("D3").Value = "3.426.26"
$ws.Range("D5").Value = "582.96"
$ws.Range("D6").Value = "134.11"
$ws.Range("D8").Value = "3.425.32"
$ws.Range("D9").Value = "0.484"
$ws.Range("D11").Value = "7.01"
$ws.Range("D13").Value = "4.013.85"
$ws.Range("D15").Value = "26.27"
$ws.Range("D16").Value = "3.427.07"
$ws.Range("D18").Value = "64.476.91"
$ws.Range("D22").Value = "380.54"
$ws.Range("D26").Value = "71.86"
$ws.Range("D27").Value = "3.568.24"
$ws.Range("D29").Value = "0.999"
$ws.Range("D31").Value = "8.07"
$ws.Range("D33").Value = "3.442.61"
$ws.Range("D35").Value = "22.96"
$ws.Range("D37").Value = "170.87"
$ws.Range("D39").Value = "6.72"
$ws.Range("D40").Value = "1.45"
$ws.Range("D42").Value = "0.0759"
$ws.Range("D43").Value = "0.800"
$ws.Range("D45").Value = "41.97"
$ws.Range("D46").Value = "4.28"
$ws.Range("D47").Value = "1.59"
$ws.Range("D48").Value = "1.10"
$ws.Range("D49").Value = "22.63"
$ws.Range("D51").Value = "2.196.07"

# Restore the default (unstyled) cell style now that the text value is set.
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).Style = "Normal"
}

# Volume(1h) (column E) values are never numeric-parseable (leading/trailing
# spaces + "%"), so they can be set directly as text.
$ws.Range("E2").Value = "  -2.98%  "
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  -3.97%  "
$ws.Range("E6").Value = "  -6.60%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -3.01%  "
$ws.Range("E9").Value = "  -6.58%  "
$ws.Range("E10").Value = "  -8.17%  "
$ws.Range("E11").Value = "  -9.15%  "
$ws.Range("E12").Value = "  -8.79%  "
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("E14").Value = "  -8.67%  "
$ws.Range("E15").Value = "  -8.56%  "
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("E18").Value = "  -2.80%  "
$ws.Range("E19").Value = "  -11.87%  "
$ws.Range("E20").Value = "  -8.57%  "
$ws.Range("E21").Value = "  -7.39%  "
$ws.Range("E22").Value = "  -10.28%  "
$ws.Range("E23").Value = "  -8.25%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("E26").Value = "  -7.04%  "
$ws.Range("E27").Value = "  -2.79%  "
$ws.Range("E28").Value = "  -8.29%  "
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("E30").Value = "  -10.07%  "
$ws.Range("E31").Value = "  -10.19%  "
$ws.Range("E32").Value = "  -11.21%  "
$ws.Range("E33").Value = "  -2.73%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  -5.58%  "
$ws.Range("E36").Value = "  -9.59%  "
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("E38").Value = "  -12.26%  "
$ws.Range("E39").Value = "  -11.29%  "
$ws.Range("E40").Value = "  -11.11%  "
$ws.Range("E41").Value = "  -10.91%  "
$ws.Range("E42").Value = "  -7.69%  "
$ws.Range("E43").Value = "  -7.00%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  -7.51%  "
$ws.Range("E46").Value = "  -14.29%  "
$ws.Range("E47").Value = "  -10.06%  "
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("E50").Value = "  -8.15%  "
$ws.Range("E51").Value = "  -5.36%  "
